# Edit script: replace sensor dataset (rows 2-5) with a new 1000-sample batch window,
# drop the now-unused 6th sample row, and custom-accuracy-widen a few columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove last data row (row 6); this shrinks the sheet dimension from AH6 to AH5.
$ws.Rows.Item(6).Delete()

# custom accuracy: widen a few numeric columns from 7 to 8 character units.
$ws.Columns("Q").ColumnWidth = 7.166666666666667
$ws.Columns("AA").ColumnWidth = 7.166666666666667
$ws.Columns("AC").ColumnWidth = 7.166666666666667

# Row 2: new batch values (Time + J1..J33 readings)
$row2 = New-Object 'object[,]' 1,34
$row2[0,0] = 45173.50694444445
$row2[0,1] = 17.776
$row2[0,2] = 11.78
$row2[0,3] = 4.002
$row2[0,4] = 37.831
$row2[0,5] = 29.945
$row2[0,6] = 13.989
$row2[0,7] = 43.56
$row2[0,8] = 21.524
$row2[0,9] = 8.901999999999999
$row2[0,10] = 13.298
$row2[0,11] = 14.882
$row2[0,12] = 15.457
$row2[0,13] = 4.465
$row2[0,14] = 13.911
$row2[0,15] = 19.348
$row2[0,16] = 12.001
$row2[0,17] = 3.404
$row2[0,18] = 2.232
$row2[0,19] = 204.374
$row2[0,20] = 38.641
$row2[0,21] = 12.84
$row2[0,22] = 25.287
$row2[0,23] = 12.826
$row2[0,24] = 3.128
$row2[0,25] = 22.227
$row2[0,26] = 11.342
$row2[0,27] = 10.33
$row2[0,28] = 12.137
$row2[0,29] = 15.411
$row2[0,30] = 3.318
$row2[0,31] = 38.832
$row2[0,32] = 6.978
$row2[0,33] = 16.052
$ws.Range("A2:AH2").Value = $row2

# Row 3: new batch values (Time + J1..J33 readings)
$row3 = New-Object 'object[,]' 1,34
$row3[0,0] = 45173.51388888889
$row3[0,1] = 16.815
$row3[0,2] = 11.911
$row3[0,3] = 1.824
$row3[0,4] = 36.401
$row3[0,5] = 29.274
$row3[0,6] = 13.233
$row3[0,7] = 50.673
$row3[0,8] = 20.36
$row3[0,9] = 8.805999999999999
$row3[0,10] = 12.939
$row3[0,11] = 14.56
$row3[0,12] = 15.229
$row3[0,13] = 4.227
$row3[0,14] = 13.159
$row3[0,15] = 18.547
$row3[0,16] = 11.34
$row3[0,17] = 1.54
$row3[0,18] = 1.092
$row3[0,19] = 192.972
$row3[0,20] = 36.819
$row3[0,21] = 12.146
$row3[0,22] = 24.411
$row3[0,23] = 12.674
$row3[0,24] = 2.374
$row3[0,25] = 24.745
$row3[0,26] = 10.728
$row3[0,27] = 9.691000000000001
$row3[0,28] = 11.371
$row3[0,29] = 15.167
$row3[0,30] = 1.247
$row3[0,31] = 46.217
$row3[0,32] = 6.7
$row3[0,33] = 15.185
$ws.Range("A3:AH3").Value = $row3

# Row 4: new batch values (Time + J1..J33 readings)
$row4 = New-Object 'object[,]' 1,34
$row4[0,0] = 45173.52083333334
$row4[0,1] = 3.363
$row4[0,2] = 2.059
$row4[0,3] = 0.885
$row4[0,4] = 7.239
$row4[0,5] = 5.375
$row4[0,6] = 2.648
$row4[0,7] = 16.021
$row4[0,8] = 4.072
$row4[0,9] = 1.685
$row4[0,10] = 2.203
$row4[0,11] = 2.896
$row4[0,12] = 2.949
$row4[0,13] = 0.86
$row4[0,14] = 2.632
$row4[0,15] = 3.686
$row4[0,16] = 2.521
$row4[0,17] = 0.931
$row4[0,18] = 0.442
$row4[0,19] = 32.808
$row4[0,20] = 7.721
$row4[0,21] = 2.429
$row4[0,22] = 4.942
$row4[0,23] = 2.43
$row4[0,24] = 0.779
$row4[0,25] = 7.205
$row4[0,26] = 2.146
$row4[0,27] = 2.087
$row4[0,28] = 2.422
$row4[0,29] = 2.958
$row4[0,30] = 0.766
$row4[0,31] = 15.09
$row4[0,32] = 1.216
$row4[0,33] = 3.041
$ws.Range("A4:AH4").Value = $row4

# Row 5: new batch values (Time + J1..J33 readings)
$row5 = New-Object 'object[,]' 1,34
$row5[0,0] = 45173.52777777778
$row5[0,1] = 5.28
$row5[0,2] = 3.62
$row5[0,3] = 0.73
$row5[0,4] = 11.45
$row5[0,5] = 8.98
$row5[0,6] = 4.16
$row5[0,7] = 16.34
$row5[0,8] = 6.4
$row5[0,9] = 2.69
$row5[0,10] = 3.88
$row5[0,11] = 4.6
$row5[0,12] = 4.78
$row5[0,13] = 1.33
$row5[0,14] = 4.14
$row5[0,15] = 5.74
$row5[0,16] = 3.7
$row5[0,17] = 0.71
$row5[0,18] = 0.39
$row5[0,19] = 55.59
$row5[0,20] = 11.51
$row5[0,21] = 3.82
$row5[0,22] = 7.49
$row5[0,23] = 3.91
$row5[0,24] = 0.87
$row5[0,25] = 7.64
$row5[0,26] = 3.37
$row5[0,27] = 3.11
$row5[0,28] = 3.64
$row5[0,29] = 4.77
$row5[0,30] = 0.55
$row5[0,31] = 14.75
$row5[0,32] = 2.04
$row5[0,33] = 4.77
$ws.Range("A5:AH5").Value = $row5

